$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.769444465637207
$ws.Range("B1").Value = 1.993281245231628
$ws.Range("C1").Value = 1.798549056053162
$ws.Range("D1").Value = 1.812283992767334
$ws.Range("E1").Value = 1.709228754043579
